$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Cell 1: NO 1 -> 29, ABU NAWAS -> EDI SATRIYO ---
$cell1 = $t.Cell(1, 1)
$cell1.Range.Find.Execute("1", $true, $true, $false, $false, $false, $true, 1, $false, "29", 2) | Out-Null
$cell1 = $t.Cell(1, 1)
$cell1.Range.Find.Execute("ABU NAWAS", $true, $true, $false, $false, $false, $true, 1, $false, "EDI SATRIYO", 2) | Out-Null

# --- Cell 2: NO 2 -> 33, ACHMAD UBAIDILLAH -> MARLION ---
$cell2 = $t.Cell(1, 2)
$cell2.Range.Find.Execute("2", $true, $true, $false, $false, $false, $true, 1, $false, "33", 2) | Out-Null
$cell2 = $t.Cell(1, 2)
$cell2.Range.Find.Execute("ACHMAD UBAIDILLAH", $true, $true, $false, $false, $false, $true, 1, $false, "MARLION", 2) | Out-Null

# --- Cell 3: NO 3 -> 35, AMRUN SAIFUDDIN -> IBRAHIM, TOPI 56 -> 58 ---
$cell3 = $t.Cell(1, 3)
$cell3.Range.Find.Execute("3", $true, $true, $false, $false, $false, $true, 1, $false, "35", 2) | Out-Null
$cell3 = $t.Cell(1, 3)
$cell3.Range.Find.Execute("AMRUN SAIFUDDIN", $true, $true, $false, $false, $false, $true, 1, $false, "IBRAHIM", 2) | Out-Null
$cell3 = $t.Cell(1, 3)
$cell3.Range.Find.Execute("56", $true, $true, $false, $false, $false, $true, 1, $false, "58", 2) | Out-Null
